# "switched back to linear models and made relevant changes"
# Rebuilds the soil_organics ANOVA/LM summary table: each response variable
# (Soil C, Soil N, Soil CN, SWR) now gets its own Df / F value / Pr(>F)
# triplet of columns instead of sharing one Df column with Chisq / Pr(>Chisq).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 — response-variable group headers.
# Old: (C1=Soil C) (E1=Soil N) (G1=Soil CN) (I1=SWR)
# New: (B1=Soil C) (E1=Soil N) (H1=Soil CN) (K1=SWR)
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Soil C"
$ws.Range("C1").ClearContents()
$ws.Range("E1").Value = "Soil N"
$ws.Range("G1").ClearContents()
$ws.Range("H1").Value = "Soil CN"
$ws.Range("I1").ClearContents()
$ws.Range("K1").Value = "SWR"

# ---------------------------------------------------------------------
# Row 2 — per-column sub headers: Df / F value / Pr(>F), repeated per
# response variable. Previously these cells carried style s="4"; in the
# new layout they go back to the default (unstyled) look.
# ---------------------------------------------------------------------
$row2 = @("B2","C2","D2","E2","F2","G2","H2","I2","J2","K2","L2","M2")
foreach ($ref in $row2) { $ws.Range($ref).ClearFormats() }

$ws.Range("B2").Value = "Df"
$ws.Range("C2").Value = "F value"
$ws.Range("D2").Value = "Pr(>F)"
$ws.Range("E2").Value = "Df"
$ws.Range("F2").Value = "F value"
$ws.Range("G2").Value = "Pr(>F)"
$ws.Range("H2").Value = "Df"
$ws.Range("I2").Value = "F value"
$ws.Range("J2").Value = "Pr(>F)"
$ws.Range("K2").Value = "Df"
$ws.Range("L2").Value = "F value"
$ws.Range("M2").Value = "Pr(>F)"

# ---------------------------------------------------------------------
# Helper data for the three model-term rows (3=Elevation, 4=Fire,
# 5=Elevation*Fire). Each response variable contributes a Df cell
# (unstyled, integer) followed by an F value / Pr(>F) pair formatted
# "0.000" (style s="5").
# ---------------------------------------------------------------------
$dfCols   = @("B","E","H","K")
$statCols = @(@("C","D"), @("F","G"), @("I","J"), @("L","M"))

$data = @{
  3 = @{ df = @(1,1,1,1);
         stats = @(
           @(6.94155590276785,   0.0137766418243072),
           @(0.35756340798521102,0.55596884905989097),
           @(5.7079483483174096, 0.0258936760329842),
           @(0.755513307733644,  0.39049465896663899)
         ) }
  4 = @{ df = @(1,1,1,1);
         stats = @(
           @(2.7182019826200401, 0.110798794125155),
           @(0.25956162041039199,0.61549409377008901),
           @(1.49258975643909,   0.23474518843320599),
           @(12.400324080444101, 0.00118586084428611)
         ) }
  5 = @{ df = @(1,1,1,1);
         stats = @(
           @(0.40421078043050002,0.53027512608080996),
           @(1.15311321231708,   0.29453587977649898),
           @(2.77139380168496,   0.110141735300711),
           @(12.980857957539101, 0.00094320950120857)
         ) }
}

foreach ($r in 3,4,5) {
  $rowData = $data[$r]
  for ($i = 0; $i -lt 4; $i++) {
    $dfRef = $dfCols[$i] + $r
    $ws.Range($dfRef).ClearFormats()
    $ws.Range($dfRef).Value = $rowData.df[$i]

    $fCol  = $statCols[$i][0] + $r
    $pCol  = $statCols[$i][1] + $r
    $ws.Range($fCol).NumberFormat = "0.000"
    $ws.Range($fCol).Value = $rowData.stats[$i][0]
    $ws.Range($pCol).NumberFormat = "0.000"
    $ws.Range($pCol).Value = $rowData.stats[$i][1]
  }
}

# ---------------------------------------------------------------------
# Row 6 — residual degrees of freedom per response variable (replaces
# the old fully-blank spacer row). The non-Df columns are now fully
# empty cells (no leftover style) rather than styled blanks.
# ---------------------------------------------------------------------
$ws.Range("B6").ClearFormats()
$ws.Range("B6").Value = 27
$ws.Range("C6").Clear() | Out-Null
$ws.Range("D6").Clear() | Out-Null
$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = 22
$ws.Range("F6").Clear() | Out-Null
$ws.Range("G6").Clear() | Out-Null
$ws.Range("H6").ClearFormats()
$ws.Range("H6").Value = 22
$ws.Range("I6").Clear() | Out-Null
$ws.Range("J6").Clear() | Out-Null
$ws.Range("K6").Value = 36

# ---------------------------------------------------------------------
# Column widths — bestFit widths recomputed by Excel for the new,
# narrower per-variable columns (was: B=3, C:J=12.1640625).
# ---------------------------------------------------------------------
$colWidths = @{
  2  = 5.6640625
  3  = 7
  4  = 6.1640625
  5  = 6
  6  = 7
  7  = 6.1640625
  8  = 7
  9  = 7
  10 = 6.1640625
  11 = 5.1640625
  12 = 7
  13 = 6.1640625
}
foreach ($col in $colWidths.Keys) {
  $ws.Columns.Item($col).ColumnWidth = $colWidths[$col] - 0.8333333333333333
}

# ---------------------------------------------------------------------
# Selection moves from the old C3:J5 block to the new L6:M6 cell.
# ---------------------------------------------------------------------
$ws.Range("L6:M6").Select() | Out-Null
